$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Folder Inventory" sheet: a new folder was scanned more recently than
#    everything else, so it becomes the new row 2 and every existing data
#    row shifts down by one (old row 2 -> new row 3, ... old row 76 -> new
#    row 77).
# ---------------------------------------------------------------------------
$wsInventory = $wb.Worksheets.Item("Folder Inventory")

$wsInventory.Rows.Item(2).Insert()
# The freshly inserted row inherits the header row's bold/border formatting
# from Excel's "insert copies format from above" behaviour; the new data
# row should look like any other plain data row, so strip that back off.
$wsInventory.Rows.Item(2).ClearFormats()

$wsInventory.Cells.Item(2, 1).Value = "Getting_started_with_Azure_AI_services"
$wsInventory.Cells.Item(2, 2).Value = "Getting_started_with_Azure_AI_services"
$wsInventory.Cells.Item(2, 3).Value = "2025-06-16 10:58:16 +0530"
$wsInventory.Cells.Item(2, 4).Value = 1
$wsInventory.Cells.Item(2, 5).Value = "Root"

# ---------------------------------------------------------------------------
# 2. "Metadata" sheet: refresh the run timestamp, folder count and workflow
#    run number.
# ---------------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Cells.Item(3, 2).Value = "2025-06-16 05:28:38 UTC"
$wsMetadata.Cells.Item(4, 2).Value = 76
# Leading apostrophe forces this to stay text (matching the source data,
# which stores the run number as a string) instead of being auto-converted
# to a number by Excel's smart entry; ClearFormats drops the resulting
# "stored as text" quote-prefix marker so the cell's formatting matches an
# untouched text cell.
$wsMetadata.Cells.Item(5, 2).Value = "'9"
$wsMetadata.Cells.Item(5, 2).ClearFormats()

# ---------------------------------------------------------------------------
# 3. "Summary" sheet: folder totals and most-recent-update now reflect the
#    new folder.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(2, 2).Value = 76
$wsSummary.Cells.Item(3, 2).Value = 76
$wsSummary.Cells.Item(5, 2).Value = "2025-06-16 10:58:16 +0530"
